$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Tarea2" card from the "In progress" column (C7) to the
#     "Review" column (D7) on the kanban board ---

# Capture the rendered text of C7 before we touch anything so the move
# reuses the existing shared string (and its rich-text runs) instead of
# creating a new, duplicate entry.
$taskText = $ws.Range("C7").Text

# Write it into the destination cell first...
$ws.Range("D7").Value = $taskText
# ...then clear the source cell now that the text lives in D7.
$ws.Range("C7").ClearContents()

# D7 becomes a proper "card" cell: centered, wrapped, same font size (11)
# the card used while it was in the "In progress" column.
$d7 = $ws.Range("D7")
$d7.Font.Size = 11
$d7.HorizontalAlignment = -4108   # xlCenter
$d7.VerticalAlignment = -4108     # xlCenter
$d7.WrapText = $true

# C7 is empty again, so it no longer needs to wrap text.
$ws.Range("C7").WrapText = $false

# Move the selection / view to the freshly edited cell.
$ws.Range("D7").Select()
